# Auto-generated Excel COM-interop script to apply the cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.953.80'
$ws.Range('E2').Value = '  +0.25%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.757.50'
$ws.Range('E3').Value = '  +0.22%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.73'
$ws.Range('E5').Value = '  -1.43%  '

$ws.Range('E6').Value = '  +0.16%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5206'
$ws.Range('E7').Value = '  +2.23%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2739'
$ws.Range('E8').Value = '  -0.90%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06165'
$ws.Range('E9').Value = '  -0.70%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.770.05'
$ws.Range('E10').Value = '  +0.94%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07035'
$ws.Range('E11').Value = '  +0.99%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.58'
$ws.Range('E12').Value = '  -0.63%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6357'
$ws.Range('E13').Value = '  +3.91%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.523'
$ws.Range('E14').Value = '  -0.15%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.74'
$ws.Range('E15').Value = '  +0.30%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  +0.07%  '

$ws.Range('E17').Value = '  +0.13%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.963.27'
$ws.Range('E18').Value = '  +0.28%  '

$ws.Range('E19').Value = '  -0.58%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006689'
$ws.Range('E20').Value = '  -3.72%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.995.94'
$ws.Range('E21').Value = '  +1.38%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.054'
$ws.Range('E22').Value = '  -0.76%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.469'
$ws.Range('E23').Value = '  +2.93%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.175'
$ws.Range('E24').Value = '  -1.87%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '138.94'
$ws.Range('E25').Value = '  +0.71%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.502'
$ws.Range('E26').Value = '  +0.63%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.841'
$ws.Range('E27').Value = '  +1.35%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.09'
$ws.Range('E28').Value = '  +0.10%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '103.31'
$ws.Range('E29').Value = '  -0.37%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08378'
$ws.Range('E30').Value = '  +2.12%  '

$ws.Range('E31').Value = '  -1.01%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.425'
$ws.Range('E32').Value = '  -2.01%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04457'
$ws.Range('E33').Value = '  -1.84%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.625'
$ws.Range('E34').Value = '  -0.57%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9931'
$ws.Range('E35').Value = '  +0.09%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6034'
$ws.Range('E36').Value = '  -1.42%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.716'
$ws.Range('E37').Value = '  +0.18%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01587'
$ws.Range('E38').Value = '  +1.74%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.954'
$ws.Range('E39').Value = '  +2.93%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.001'
$ws.Range('E40').Value = '  +0.09%  '

$ws.Range('B41').Value = 'PaxosStandard'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.002'
$ws.Range('E41').Value = '  +0.07%  '

$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '102.39'
$ws.Range('E42').Value = '  -1.60%  '

$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3861'
$ws.Range('E43').Value = '  -0.64%  '

$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7404'
$ws.Range('E44').Value = '  -0.26%  '

$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.913'
$ws.Range('E45').Value = '  -0.48%  '

$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05507'
$ws.Range('E46').Value = '  +1.51%  '

$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.288'
$ws.Range('E47').Value = '  +4.76%  '

$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1111'
$ws.Range('E48').Value = '  -0.29%  '

$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.09'
$ws.Range('E49').Value = '  -0.07%  '

$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.39'
$ws.Range('E50').Value = '  -0.79%  '

$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.003'
$ws.Range('E51').Value = '  +0.58%  '
